$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename "_old" header columns (A1:J1) to "_FV2210" ---
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
for ($i = 0; $i -lt $fv2210Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}

# --- 2. Rename "_new" header columns (L1:U1) to "_FV2304" (K1 "diff" stays the same) ---
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $fv2304Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# --- 3. Turn the used range into an Excel Table ("Table1") so the renamed
#        headers become the table's column headers ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U89"), $null, 1)
$tbl.Name = "Table1"

# --- 4. Freeze the header row (pane split below row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Host "edit complete"
